$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 131.08333
$ws.Range("I9").Value = 96.28570999999999
$ws.Range("J9").Value = 179.8
$ws.Range("K9").Value = 96.28570999999999
$ws.Range("L9").Value = 179.8
$ws.Range("M9").Value = 72.71429000000001
$ws.Range("N9").Value = -517.8
$ws.Range("H80").Value = 11230796
$ws.Range("I80").Value = 11111396
$ws.Range("J80").Value = 11365120
$ws.Range("K80").Value = 33334188
$ws.Range("L80").Value = 34095360
$ws.Range("M80").Value = -33333190
$ws.Range("N80").Value = -34097356
$ws.Range("H83").Value = 11230796
$ws.Range("I83").Value = 11111396
$ws.Range("J83").Value = 11365120
$ws.Range("K83").Value = 100002564
$ws.Range("L83").Value = 102286080
$ws.Range("M83").Value = -99997572
$ws.Range("N83").Value = -102296064
$ws.Range("H132").Value = 173204.39
$ws.Range("I132").Value = 3624.2886
$ws.Range("J132").Value = 1432942.2
$ws.Range("K132").Value = 10872.8658
$ws.Range("L132").Value = 4298826.6
$ws.Range("M132").Value = -8342.8658
$ws.Range("N132").Value = -4303886.6
$ws.Range("H138").Value = 131846.61
$ws.Range("I138").Value = 2332.6667
$ws.Range("J138").Value = 162683.25
$ws.Range("K138").Value = 6998.000100000001
$ws.Range("L138").Value = 488049.75
$ws.Range("M138").Value = -1858.000100000001
$ws.Range("N138").Value = -498329.75

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 25584.072
$ws.Range("I32").Value = 19454.365
$ws.Range("J32").Value = 43535.355
$ws.Range("K32").Value = 19454.365
$ws.Range("L32").Value = 43535.355
$ws.Range("M32").Value = -19167.365
$ws.Range("N32").Value = -44109.355
$ws.Range("H61").Value = 2079.4
$ws.Range("I61").Value = 1365.2
$ws.Range("J61").Value = 3507.8
$ws.Range("K61").Value = 1365.2
$ws.Range("L61").Value = 3507.8
$ws.Range("M61").Value = -1153.2
$ws.Range("N61").Value = -3931.8
$ws.Range("J63").Value = 3000
$ws.Range("L63").Value = 3000
$ws.Range("N63").Value = -4372
$ws.Range("J66").Value = 3000
$ws.Range("L66").Value = 15000
$ws.Range("N66").Value = -21864
$ws.Range("H102").Value = 2273.3333
$ws.Range("I102").Value = 820
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 820
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = 802
$ws.Range("N102").Value = -6244
$ws.Range("H136").Value = 2079.4
$ws.Range("I136").Value = 1365.2
$ws.Range("J136").Value = 3507.8
$ws.Range("K136").Value = 4095.6
$ws.Range("L136").Value = 10523.4
$ws.Range("M136").Value = -1545.6
$ws.Range("N136").Value = -15623.4

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1709.909
$ws.Range("I20").Value = 1501.5
$ws.Range("J20").Value = 1960
$ws.Range("K20").Value = 1501.5
$ws.Range("L20").Value = 1960
$ws.Range("M20").Value = -1254.5
$ws.Range("N20").Value = -2454
$ws.Range("H94").Value = 3592.9
$ws.Range("I94").Value = 627.53845
$ws.Range("J94").Value = 9100
$ws.Range("K94").Value = 627.53845
$ws.Range("L94").Value = 9100
$ws.Range("M94").Value = -176.53845
$ws.Range("N94").Value = -10002
$ws.Range("H99").Value = 2500
$ws.Range("J99").Value = 2800
$ws.Range("L99").Value = 2800
$ws.Range("N99").Value = -5796

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 415.4
$ws.Range("I22").Value = 338.30768
$ws.Range("J22").Value = 558.5714
$ws.Range("K22").Value = 338.30768
$ws.Range("L22").Value = 558.5714
$ws.Range("M22").Value = 11.69232
$ws.Range("N22").Value = -1258.5714
$ws.Range("H31").Value = 31702.863
$ws.Range("I31").Value = 42833.76
$ws.Range("J31").Value = 17056.947
$ws.Range("K31").Value = 42833.76
$ws.Range("L31").Value = 17056.947
$ws.Range("M31").Value = -42538.76
$ws.Range("N31").Value = -17646.947
$ws.Range("H34").Value = 31702.863
$ws.Range("I34").Value = 42833.76
$ws.Range("J34").Value = 17056.947
$ws.Range("K34").Value = 42833.76
$ws.Range("L34").Value = 17056.947
$ws.Range("M34").Value = -42631.76
$ws.Range("N34").Value = -17460.947

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 616.5454999999999
$ws.Range("I34").Value = 80
$ws.Range("J34").Value = 735.7778
$ws.Range("K34").Value = 240
$ws.Range("L34").Value = 2207.3334
$ws.Range("M34").Value = -156
$ws.Range("N34").Value = -2375.3334
$ws.Range("H39").Value = 35069
$ws.Range("J39").Value = 35069
$ws.Range("L39").Value = 105207
$ws.Range("N39").Value = -105795
$ws.Range("H55").Value = 122224390
$ws.Range("J55").Value = 122224390
$ws.Range("L55").Value = 366673170
$ws.Range("N55").Value = -366673524
$ws.Range("H113").Value = 526.2692
$ws.Range("I113").Value = 461.53845
$ws.Range("J113").Value = 591
$ws.Range("K113").Value = 1384.61535
$ws.Range("L113").Value = 1773
$ws.Range("M113").Value = 785.38465
$ws.Range("N113").Value = -6113
$ws.Range("H131").Value = 189555.81
$ws.Range("J131").Value = 204992.4
$ws.Range("L131").Value = 614977.2
$ws.Range("N131").Value = -625057.2

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 640.9167
$ws.Range("I16").Value = 653.7273
$ws.Range("J16").Value = 500
$ws.Range("K16").Value = 653.7273
$ws.Range("L16").Value = 500
$ws.Range("M16").Value = -483.7273
$ws.Range("N16").Value = -840
$ws.Range("H22").Value = 841.7
$ws.Range("I22").Value = 755.8
$ws.Range("J22").Value = 927.6
$ws.Range("K22").Value = 755.8
$ws.Range("L22").Value = 927.6
$ws.Range("M22").Value = -460.8
$ws.Range("N22").Value = -1517.6
$ws.Range("H27").Value = 841.7
$ws.Range("I27").Value = 755.8
$ws.Range("J27").Value = 927.6
$ws.Range("K27").Value = 755.8
$ws.Range("L27").Value = 927.6
$ws.Range("M27").Value = -648.8
$ws.Range("N27").Value = -1141.6
$ws.Range("H82").Value = 1491.4572
$ws.Range("I82").Value = 1215.0476
$ws.Range("J82").Value = 1906.0714
$ws.Range("K82").Value = 1215.0476
$ws.Range("L82").Value = 1906.0714
$ws.Range("M82").Value = -854.0476000000001
$ws.Range("N82").Value = -2628.0714
$ws.Range("H85").Value = 1491.4572
$ws.Range("I85").Value = 1215.0476
$ws.Range("J85").Value = 1906.0714
$ws.Range("K85").Value = 1215.0476
$ws.Range("L85").Value = 1906.0714
$ws.Range("M85").Value = 32.9523999999999
$ws.Range("N85").Value = -4402.0714

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H127").Value = 50000
$ws.Range("J127").Value = 50000
$ws.Range("L127").Value = 50000
$ws.Range("N127").Value = -59920
$ws.Range("H132").Value = 2859.8394
$ws.Range("I132").Value = 664.06525
$ws.Range("J132").Value = 12960.4
$ws.Range("K132").Value = 1992.19575
$ws.Range("L132").Value = 38881.2
$ws.Range("M132").Value = 537.8042500000001
$ws.Range("N132").Value = -43941.2
